# Issue #30: Add option to not convert text to numbers: convertTextToNumber
#
# Adds a new worksheet "28" (for the new regression test data) right after
# the existing "23" sheet, populates its header row, and leaves it as the
# active/selected sheet - matching the authored commit.

$wb = $excel.ActiveWorkbook

# The existing (only) worksheet, "23".
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet right after "23" and rename it "28".
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "28"

# Populate header row with the new shared-string values.
$ws2.Cells.Item(1, 1).Value = "furtherInformation.icon"
$ws2.Cells.Item(1, 2).Value = "furtherInformation.description"
$ws2.Cells.Item(1, 3).Value = "button.title"
$ws2.Cells.Item(1, 4).Value = "button.link"

# Match the authored selection on the new (now active) sheet.
$ws2.Range("D3").Select() | Out-Null
